$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "RM 232" row (row 26) and the "SC 92" row (row 28) were removed from the
# dataset entirely, shifting all subsequent rows up. Delete the lower-indexed
# row second so the earlier row number (26) is still valid when we get to it.
$ws.Rows(28).Delete()
$ws.Rows(26).Delete()

# After the two rows were removed, several cells in column F (the re-imputed
# "missing" column) were recomputed: some previously-missing cells now carry
# a value, and some previously-filled cells are now missing again.
$ws.Range("F2").Value = 18.03
$ws.Range("F6").Value = ""
$ws.Range("F12").Value = 17.45
$ws.Range("F14").Value = ""
$ws.Range("F20").Value = 17.73
$ws.Range("F21").Value = 16.58
$ws.Range("F22").Value = ""
$ws.Range("F23").Value = ""
$ws.Range("F31").Value = 17.18
$ws.Range("F33").Value = 17.53

# Two rows also had their column E (re-)imputed value swapped.
$ws.Range("E30").Value = -5.7
$ws.Range("E32").Value = ""
